$wb = $excel.ActiveWorkbook

# ===== Sheet: 展览 =====
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value = 115
$ws.Range("F3").Value = 1012
$ws.Range("F6").Value = 1572
$ws.Range("G6").Value = "已售罄"
$ws.Range("F7").Value = 40489
$ws.Range("G7").Value = "已售罄"
$ws.Range("F10").Value = 8807
$ws.Range("F11").Value = 180
$ws.Range("F12").Value = 615
$ws.Range("F13").Value = 766
$ws.Range("F14").Value = 621
$ws.Range("F15").Value = 143
$ws.Range("F16").Value = 238
$ws.Range("F17").Value = 746
$ws.Range("F20").Value = 632
$ws.Range("F21").Value = 261
$ws.Range("F22").Value = 1189
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 453
$ws.Range("F25").Value = 615
$ws.Range("F26").Value = 414
$ws.Range("F27").Value = 614
$ws.Range("F28").Value = 641
$ws.Range("F30").Value = 33
$ws.Range("F31").Value = 38
$ws.Range("F33").Value = 434
$ws.Range("F34").Value = 20
$ws.Range("F35").Value = 174
$ws.Range("F36").Value = 859
$ws.Range("F37").Value = 392
$ws.Range("F38").Value = 43
$ws.Range("F39").Value = 192
$ws.Range("F40").Value = 86
$ws.Range("F41").Value = 288
$ws.Range("F42").Value = 1097
$ws.Range("F43").Value = 234
$ws.Range("F44").Value = 1115
$ws.Range("F45").Value = 353
$ws.Range("F46").Value = 85
$ws.Range("F47").Value = 15
$ws.Range("F48").Value = 25
$ws.Range("F49").Value = 51

# ===== Sheet: 演出 =====
$ws = $wb.Worksheets.Item("演出")

$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 20
$ws.Range("F4").Value = 203
$ws.Range("F5").Value = 338
$ws.Range("F6").Value = 4407
$ws.Range("F7").Value = 9
$ws.Range("F20").Value = 4373

# ===== Sheet: 本地生活 =====
$ws = $wb.Worksheets.Item("本地生活")

$ws.Range("F2").Value = 1888
$ws.Range("F3").Value = 442
$ws.Range("F4").Value = 424
$ws.Range("F5").Value = 228

# ===== Sheet: 全部类型 =====
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F2").Value = 1888
$ws.Range("F3").Value = 442
$ws.Range("F4").Value = 424
$ws.Range("F5").Value = 115
$ws.Range("F6").Value = 1012
$ws.Range("F7").Value = 5
$ws.Range("B8").Value = "2024-07-20"
$ws.Range("C8").Value = "广州·KKWORLD-【陈张太康】配音演员签名内场礼包"
$ws.Range("D8").Value = "新港东路1000号 保利世贸博览馆"
$ws.Range("E8").Value = "2024.07.20 10:30-07.20 14:00"
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 298
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89072"
$ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202407/CTKJTCjG1720513282520.png"
$ws.Range("B9").Value = "2024-07-20"
$ws.Range("C9").Value = "广州·冰兔2024线下live「过去和未来」"
$ws.Range("D9").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws.Range("E9").Value = "2024.07.20 20:00-07.20 22:00"
$ws.Range("F9").Value = 203
$ws.Range("G9").Value = 198
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"
$ws.Range("C10").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws.Range("D10").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws.Range("E10").Value = "2024.07.20 19:30-07.20 21:10"
$ws.Range("F10").Value = 338
$ws.Range("G10").Value = 480
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"
$ws.Range("B11").Value = "2024-07-21"
$ws.Range("C11").Value = "广州·昨日重现——唯美英文经典歌曲演唱会"
$ws.Range("D11").Value = "东风中路299号 广州中山纪念堂"
$ws.Range("E11").Value = "2024.07.21 19:30-07.21 21:30"
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=86802"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"
$ws.Range("C12").Value = "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws.Range("D12").Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws.Range("E12").Value = "2024.07.21 14:30-07.21 16:00"
$ws.Range("F12").Value = 317
$ws.Range("G12").Value = 280
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=87034"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"
$ws.Range("F13").Value = 8807
$ws.Range("F14").Value = 180
$ws.Range("F15").Value = 615
$ws.Range("F17").Value = 228
$ws.Range("F18").Value = 766
$ws.Range("F19").Value = 621
$ws.Range("F21").Value = 143
$ws.Range("F22").Value = 238
$ws.Range("F23").Value = 746
$ws.Range("F26").Value = 261
$ws.Range("F27").Value = 1189
$ws.Range("F28").Value = 4
$ws.Range("F29").Value = 453
$ws.Range("F30").Value = 414
$ws.Range("F31").Value = 614
$ws.Range("F32").Value = 641
$ws.Range("F33").Value = 33
$ws.Range("F34").Value = 38
$ws.Range("F37").Value = 434
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 174
$ws.Range("F40").Value = 859
$ws.Range("F41").Value = 392
$ws.Range("F42").Value = 43
$ws.Range("F43").Value = 192
$ws.Range("F44").Value = 86
$ws.Range("F45").Value = 288
$ws.Range("F46").Value = 234
$ws.Range("F47").Value = 1115
$ws.Range("F48").Value = 353
$ws.Range("F49").Value = 85
